# Weekly update to the "Hortaliza, Femacal de La Calera - Zapallo italiano" price table:
# re-shuffled existing rows 353-376 and appended two new price rows (377-378).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 120, 5000, 5500, 5250, "$/caja 36 unidades", "Provincia de Quillota", 146, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 110, 9500, 10000, 9727, "$/caja 70 unidades", "Provincia de Quillota", 139, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44312, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 105, 6000, 6500, 6238, "$/caja 70 unidades", "Provincia de Quillota", 89, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44399, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 230, 8000, 9000, 8478, "$/caja 70 unidades", "Región de Arica y Parinacota", 121, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44522, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 130, 4000, 4500, 4231, "$/caja 36 unidades", "Provincia de Quillota", 118, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44543, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 195, 7500, 8000, 7751, "$/caja 60 unidades", "Limache", 129, 60, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44167, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 60, 4000, 4000, 4000, "$/caja 36 unidades", "Limache", 111, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44167, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 60, 8000, 8000, 8000, "$/caja 70 unidades", "Limache", 114, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44277, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 140, 9000, 9000, 9000, "$/caja 70 unidades", "Provincia de Quillota", 129, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44258, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 95, 9000, 9000, 9000, "$/caja 70 unidades", "Provincia de Quillota", 129, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44390, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 45, 12000, 12000, 12000, "$/caja 70 unidades", "Región de Arica y Parinacota", 171, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44349, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 108, 12000, 13000, 12324, "$/caja 70 unidades", "Región de Arica y Parinacota", 176, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44285, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 110, 8500, 9000, 8773, "$/caja 70 unidades", "Provincia de Quillota", 125, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44498, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 125, 7000, 7500, 7240, "$/caja 36 unidades", "Limache", 201, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44498, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 105, 10000, 11000, 10476, "$/caja 70 unidades", "Región de Arica y Parinacota", 150, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44179, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 76, 4500, 4500, 4500, "$/caja 36 unidades", "Provincia de Quillota", 125, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44179, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 75, 8000, 8000, 8000, "$/caja 70 unidades", "Provincia de Quillota", 114, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44418, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 200, 6300, 7000, 6605, "$/caja 70 unidades", "Región de Arica y Parinacota", 94, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44595, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 160, 4500, 5000, 4750, "$/caja 36 unidades", "Provincia de Quillota", 132, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44595, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 160, 9000, 9500, 9250, "$/caja 70 unidades", "Provincia de Quillota", 132, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44335, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 190, 7000, 7800, 7421, "$/caja 70 unidades", "Región de Arica y Parinacota", 106, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44552, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 185, 4000, 4500, 4257, "$/caja 36 unidades", "Provincia de Quillota", 118, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44552, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 290, 6000, 7500, 6810, "$/caja 70 unidades", "Región de Arica y Parinacota", 97, 70, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44544, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 110, 4000, 4500, 4227, "$/caja 36 unidades", "Provincia de Quillota", 117, 36, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44544, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 230, 7000, 7800, 7452, "$/caja 60 unidades", "Provincia de Quillota", 124, 60, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44160, 5, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 68, 7500, 7500, 7500, "$/caja 70 unidades", "Región de Arica y Parinacota", 107, 70, "Hortaliza")
)

$startRow = 353
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rowsData[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

# The two brand-new rows need the Fecha (column D) date number format re-applied,
# since freshly-created cells default to General formatting.
$ws.Cells.Item(377, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(378, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
